# Update effort hours on the single Effort sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2        # Purpose, scope, definitions: 1.5 -> 2
$ws.Range("B5").Value = 3        # Product functions: 2 -> 3
$ws.Range("B6").Value = 1.15     # Domain assumptions: 0.75 -> 1.15
$ws.Range("B8").Value = 4.5      # Functional requirements: 2.5 -> 4.5

# Leave the cursor where the author left it before saving.
$ws.Range("B14").Select()
